$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddr, $text) {
    $rng = $ws.Range($rangeAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "34.487.35"
$ws.Range("E2").Value = "  -3.05%  "

# Row 3 - Ethereum
Set-TextValue "D3" "1.800.14"
$ws.Range("E3").Value = "  -2.26%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.52%  "

# Row 5 - BNB
Set-TextValue "D5" "228.52"
$ws.Range("E5").Value = "  -1.47%  "

# Row 6 - XRP
Set-TextValue "D6" "0.609"
$ws.Range("E6").Value = "  -1.56%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.57%  "

# Row 8 - Solana
Set-TextValue "D8" "38.76"
$ws.Range("E8").Value = "  -11.41%  "

# Row 9 - Cardano
Set-TextValue "D9" "0.321"
$ws.Range("E9").Value = "  +3.08%  "

# Row 10 - Dogecoin
Set-TextValue "D10" "0.0676"
$ws.Range("E10").Value = "  -3.79%  "

# Row 11 - TRON
Set-TextValue "D11" "0.0987"
$ws.Range("E11").Value = "  -2.14%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("E12").Value = "  -2.21%  "

# Row 13 - Chainlink
Set-TextValue "D13" "11.10"
$ws.Range("E13").Value = "  -1.93%  "

# Row 14 - now Polygon (was WrappedEther)
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue "D14" "0.656"
$ws.Range("E14").Value = "  -2.83%  "

# Row 15 - now WrappedEther (was Polygon)
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D15" "1.773.15"
$ws.Range("E15").Value = "  -4.10%  "

# Row 16 - Polkadot
Set-TextValue "D16" "4.54"
$ws.Range("E16").Value = "  -4.26%  "

# Row 17 - WrappedBTC
Set-TextValue "D17" "34.543.46"
$ws.Range("E17").Value = "  -2.83%  "

# Row 18 - Litecoin
Set-TextValue "D18" "68.91"
$ws.Range("E18").Value = "  -2.17%  "

# Row 19 - BitcoinCash
Set-TextValue "D19" "239.98"
$ws.Range("E19").Value = "  -2.07%  "

# Row 20 - ShibaInu
$ws.Range("E20").Value = "  -3.18%  "

# Row 21 - Avalanche
$ws.Range("E21").Value = "  -2.56%  "

# Row 22 - Uniswap
Set-TextValue "D22" "4.66"
$ws.Range("E22").Value = "  +0.39%  "

# Row 23 - Dai
$ws.Range("E23").Value = "  +0.52%  "

# Row 24 - Toncoin
Set-TextValue "D24" "2.22"
$ws.Range("E24").Value = "  -0.07%  "

# Row 25 - Monero
Set-TextValue "D25" "171.84"
$ws.Range("E25").Value = "  -0.26%  "

# Row 26 - Cosmos
Set-TextValue "D26" "7.69"

# Row 27 - EthereumClassic
Set-TextValue "D27" "17.10"
$ws.Range("E27").Value = "  -4.19%  "

# Row 28 - Stellar
$ws.Range("E28").Value = "  -1.53%  "

# Row 29 - PancakeSwap
$ws.Range("E29").Value = "  -4.81%  "

# Row 30 - BinanceUSD
$ws.Range("E30").Value = "  +0.52%  "

# Row 31 - Filecoin
Set-TextValue "D31" "4.02"
$ws.Range("E31").Value = "  +2.05%  "

# Row 32 - Hedera
Set-TextValue "D32" "0.0538"
$ws.Range("E32").Value = "  -2.66%  "

# Row 33 - InternetComputer(DFINITY)
Set-TextValue "D33" "3.86"
$ws.Range("E33").Value = "  -5.59%  "

# Row 34 - TrustWalletToken
Set-TextValue "D34" "1.23"
$ws.Range("E34").Value = "  +7.88%  "

# Row 35 - LidoDAOToken
$ws.Range("E35").Value = "  -3.91%  "

# Row 36 - ImmutableX
$ws.Range("E36").Value = "  -0.25%  "

# Row 37 - Aave
Set-TextValue "D37" "90.78"
$ws.Range("E37").Value = "  -5.53%  "

# Row 38 - WEMIXToken
Set-TextValue "D38" "1.32"
$ws.Range("E38").Value = "  +4.31%  "

# Row 39 - Maker
Set-TextValue "D39" "1.313.54"
$ws.Range("E39").Value = "  -2.62%  "

# Row 40 - VeChain
$ws.Range("E40").Value = "  -2.60%  "

# Row 41 - HuobiToken
Set-TextValue "D41" "2.43"
$ws.Range("E41").Value = "  -0.96%  "

# Row 42 - ARBITRUM
$ws.Range("E42").Value = "  -6.45%  "

# Row 43 - InjectiveProtocol
Set-TextValue "D43" "14.22"
$ws.Range("E43").Value = "  -8.06%  "

# Row 44 - RenderToken
$ws.Range("E44").Value = "  -11.50%  "

# Row 45 - MXToken
$ws.Range("E45").Value = "  -4.19%  "

# Row 46 - FraxShare
Set-TextValue "D46" "6.18"
$ws.Range("E46").Value = "  -1.89%  "

# Row 47 - Kaspa
Set-TextValue "D47" "0.0513"
$ws.Range("E47").Value = "  -1.10%  "

# Row 48 - RocketPoolETH
$ws.Range("E48").Value = "  -1.43%  "

# Row 49 - PaxDollar
$ws.Range("E49").Value = "  +0.56%  "

# Row 50 - Cronos
$ws.Range("E50").Value = "  +4.10%  "

# Row 51 - Quant
Set-TextValue "D51" "97.28"
$ws.Range("E51").Value = "  -5.22%  "
